$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated NATMI TPM values (Inhbb-Acvr1) per commit "update scripts wuth new tpm"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.372558333333333
$ws.Range("H2").Value = 4.117675
$ws.Range("I2").Value = 0.3340102211301095
$ws.Range("J2").Value = 0.3340102211301095
$ws.Range("M2").Value = 10.92359866666667
$ws.Range("N2").Value = 32.770796
$ws.Range("O2").Value = 0.2236009040380497
$ws.Range("P2").Value = 0.2236009040380497
$ws.Range("Q2").Value = 14.99327637992223
$ws.Range("R2").Value = 134.9394874193
$ws.Range("S2").Value = 0.07468498740264137
$ws.Range("T2").Value = 0.07468498740264137
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.372558333333333
$ws.Range("H3").Value = 4.117675
$ws.Range("I3").Value = 0.3340102211301095
$ws.Range("J3").Value = 0.3340102211301095
$ws.Range("O3").Value = 0.4261214970992155
$ws.Range("P3").Value = 0.4261214970992155
$ws.Range("Q3").Value = 28.57303911592222
$ws.Range("R3").Value = 257.1573520433
$ws.Range("S3").Value = 0.1423289354744023
$ws.Range("T3").Value = 0.1423289354744023
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.372558333333333
$ws.Range("H4").Value = 4.117675
$ws.Range("I4").Value = 0.3340102211301095
$ws.Range("J4").Value = 0.3340102211301095
$ws.Range("M4").Value = 13.06524766666667
$ws.Range("N4").Value = 39.195743
$ws.Range("O4").Value = 0.2674394472823625
$ws.Range("P4").Value = 0.2674394472823625
$ws.Range("Q4").Value = 17.93281456194722
$ws.Range("R4").Value = 161.395331057525
$ws.Range("S4").Value = 0.08932750892569617
$ws.Range("T4").Value = 0.08932750892569617
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.372558333333333
$ws.Range("H5").Value = 4.117675
$ws.Range("I5").Value = 0.3340102211301095
$ws.Range("J5").Value = 0.3340102211301095
$ws.Range("M5").Value = 4.046901
$ws.Range("N5").Value = 12.140703
$ws.Range("O5").Value = 0.0828381515803724
$ws.Range("P5").Value = 0.0828381515803724
$ws.Range("Q5").Value = 5.554607691725
$ws.Range("R5").Value = 49.991469225525
$ws.Range("S5").Value = 0.02766878932736972
$ws.Range("T5").Value = 0.02766878932736972
$ws.Range("I6").Value = 0.01293592767872722
$ws.Range("J6").Value = 0.01293592767872721
$ws.Range("M6").Value = 10.92359866666667
$ws.Range("N6").Value = 32.770796
$ws.Range("O6").Value = 0.2236009040380497
$ws.Range("P6").Value = 0.2236009040380497
$ws.Range("Q6").Value = 0.5806766579226668
$ws.Range("R6").Value = 5.226089921304001
$ws.Range("S6").Value = 0.002892485123534235
$ws.Range("T6").Value = 0.002892485123534235
$ws.Range("I7").Value = 0.01293592767872722
$ws.Range("J7").Value = 0.01293592767872721
$ws.Range("O7").Value = 0.4261214970992155
$ws.Range("P7").Value = 0.4261214970992155
$ws.Range("S7").Value = 0.005512276868826421
$ws.Range("T7").Value = 0.005512276868826421
$ws.Range("I8").Value = 0.01293592767872722
$ws.Range("J8").Value = 0.01293592767872721
$ws.Range("M8").Value = 13.06524766666667
$ws.Range("N8").Value = 39.195743
$ws.Range("O8").Value = 0.2674394472823625
$ws.Range("P8").Value = 0.2674394472823625
$ws.Range("Q8").Value = 0.6945224354646667
$ws.Range("R8").Value = 6.250701919182
$ws.Range("S8").Value = 0.003459577348483422
$ws.Range("T8").Value = 0.003459577348483421
$ws.Range("I9").Value = 0.01293592767872722
$ws.Range("J9").Value = 0.01293592767872721
$ws.Range("M9").Value = 4.046901
$ws.Range("N9").Value = 12.140703
$ws.Range("O9").Value = 0.0828381515803724
$ws.Range("P9").Value = 0.0828381515803724
$ws.Range("Q9").Value = 0.215125163358
$ws.Range("R9").Value = 1.936126470222
$ws.Range("S9").Value = 0.00107158833788314
$ws.Range("T9").Value = 0.00107158833788314
$ws.Range("G10").Value = 2.683614
$ws.Range("H10").Value = 8.050841999999999
$ws.Range("I10").Value = 0.6530538511911632
$ws.Range("J10").Value = 0.6530538511911632
$ws.Range("M10").Value = 10.92359866666667
$ws.Range("N10").Value = 32.770796
$ws.Range("O10").Value = 0.2236009040380497
$ws.Range("P10").Value = 0.2236009040380497
$ws.Range("Q10").Value = 29.314722312248
$ws.Range("R10").Value = 263.832500810232
$ws.Range("S10").Value = 0.1460234315118741
$ws.Range("T10").Value = 0.1460234315118741
$ws.Range("G11").Value = 2.683614
$ws.Range("H11").Value = 8.050841999999999
$ws.Range("I11").Value = 0.6530538511911632
$ws.Range("J11").Value = 0.6530538511911632
$ws.Range("O11").Value = 0.4261214970992155
$ws.Range("P11").Value = 0.4261214970992155
$ws.Range("Q11").Value = 55.865755160888
$ws.Range("R11").Value = 502.7917964479919
$ws.Range("S11").Value = 0.2782802847559868
$ws.Range("T11").Value = 0.2782802847559868
$ws.Range("G12").Value = 2.683614
$ws.Range("H12").Value = 8.050841999999999
$ws.Range("I12").Value = 0.6530538511911632
$ws.Range("J12").Value = 0.6530538511911632
$ws.Range("M12").Value = 13.06524766666667
$ws.Range("N12").Value = 39.195743
$ws.Range("O12").Value = 0.2674394472823625
$ws.Range("P12").Value = 0.2674394472823625
$ws.Range("Q12").Value = 35.062081551734
$ws.Range("R12").Value = 315.558733965606
$ws.Range("S12").Value = 0.1746523610081829
$ws.Range("T12").Value = 0.1746523610081829
$ws.Range("G13").Value = 2.683614
$ws.Range("H13").Value = 8.050841999999999
$ws.Range("I13").Value = 0.6530538511911632
$ws.Range("J13").Value = 0.6530538511911632
$ws.Range("M13").Value = 4.046901
$ws.Range("N13").Value = 12.140703
$ws.Range("O13").Value = 0.0828381515803724
$ws.Range("P13").Value = 0.0828381515803724
$ws.Range("Q13").Value = 10.860320180214
$ws.Range("R13").Value = 97.74288162192599
$ws.Range("S13").Value = 0.05409777391511954
$ws.Range("T13").Value = 0.05409777391511954
